# "Generate Report for handback"
# Updates the localization-status workbook to reflect that both files
# (ae38855b... and fd52c5fc...) have been handed back and are in sync
# with en-US, for both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet - mirrors the same "Ready for handoff" text in its
# per-language status columns (B = zh-cn, C = de-de), so it must be
# updated in lock-step with the shared string text below.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $status
$overview.Range("C2").Value = $status
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# ---------------------------------------------------------------
# zh-cn sheet ("Status" table in xl/worksheets/sheet2.xml)
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column (B) moves from "Ready for handoff" to "Handed back..."
$zh.Range("B2").Value = $status
$zh.Range("B3").Value = $status

# New "Latest Target File" (E) / "Latest Handback File" (F) columns,
# mirroring the existing source (A) / handoff (C) file names.
$zh.Range("E2").Value = "ae38855b-7f4f-4934-8340-95b809897df5.md"
$zh.Range("F2").Value = "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf"

$zh.Range("E3").Value = "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md"
$zh.Range("F3").Value = "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.zh-cn.xlf"

# Latest Handback DateTime (G) now has real timestamps instead of the
# "never handed back" placeholder.
$zh.Range("G2").Value = "2016-01-20 07:46:22"
$zh.Range("G3").Value = "2016-01-20 07:46:22"

# Rebuild the hyperlinks in row order (A2, C2, E2, F2, A3, C3, E3, F3, A4)
# so relationship ids line up the way Excel would naturally emit them.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8f100db4e9ff0febb5a5e30cf16a52568f7733d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8f100db4e9ff0febb5a5e30cf16a52568f7733d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8f100db4e9ff0febb5a5e30cf16a52568f7733d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8f100db4e9ff0febb5a5e30cf16a52568f7733d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# Give the new cells the same "hyperlink" look as the other link cells.
$zh.Range("E2,F2,E3,F3").Font.Underline = 2
$zh.Range("E2,F2,E3,F3").Font.Color = 15570276

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $status
$de.Range("B3").Value = $status

$de.Range("E2").Value = "ae38855b-7f4f-4934-8340-95b809897df5.md"
$de.Range("F2").Value = "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf"

$de.Range("E3").Value = "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md"
$de.Range("F3").Value = "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.de-de.xlf"

$de.Range("G2").Value = "2016-01-20 07:46:44"
$de.Range("G3").Value = "2016-01-20 07:46:44"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee3e112f5321d43c0f546a6efd6d073234a1469/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/ae38855b-7f4f-4934-8340-95b809897df5.md", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee3e112f5321d43c0f546a6efd6d073234a1469/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf", [Type]::Missing, [Type]::Missing, "ae38855b-7f4f-4934-8340-95b809897df5.7ca250d313629c64a369540b17ad53292a8c108d.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee3e112f5321d43c0f546a6efd6d073234a1469/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.de-de.xlf", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/e2e/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee3e112f5321d43c0f546a6efd6d073234a1469/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.de-de.xlf", [Type]::Missing, [Type]::Missing, "fd52c5fc-8869-4aa3-81d4-fe0e6091bfc5.ee01159c09dc3c768cf766dcc8fe8f5cb41f7209.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cdc19e22ffe06624a06367a8cf920888b6c0ee31/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

$de.Range("E2,F2,E3,F3").Font.Underline = 2
$de.Range("E2,F2,E3,F3").Font.Color = 15570276

Write-Host "Handback report generated."
